$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 36 (weekly update: new latest record pushed in,
# all subsequent rows shift down by one).
$ws.Rows(36).Insert()

$ws.Cells.Item(36, 1).Value = 5
$ws.Cells.Item(36, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(36, 3).Value = "Maule"
$ws.Cells.Item(36, 4).Value = 44510
$ws.Cells.Item(36, 5).Value = 7
$ws.Cells.Item(36, 6).Value = 100112022
$ws.Cells.Item(36, 7).Value = "Arveja Verde"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 400
$ws.Cells.Item(36, 11).Value = 14000
$ws.Cells.Item(36, 12).Value = 14000
$ws.Cells.Item(36, 13).Value = 14000
$ws.Cells.Item(36, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(36, 15).Value = "Región del Maule"
$ws.Cells.Item(36, 16).Value = 560
$ws.Cells.Item(36, 17).Value = 25
$ws.Cells.Item(36, 18).Value = "Hortaliza"
